$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Result")

$data = @(
    @("testMultipleLogin", "PASS", "01_03_2017_01_08_52", "CHROME"),
    @("testMultipleLogin", "PASS", "01_03_2017_01_08_59", "CHROME"),
    @("testLoginLogout", "PASS", "01_03_2017_01_09_14", "CHROME")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
